$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "65.621.07"
$ws.Range("E2").Value = "  -1.60%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.718.86"
$ws.Range("E3").Value = "  +4.15%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.08%  "

# Row 5 - now BNB (was Solana)
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").Value = "413.05"
$ws.Range("E5").Value = "  -1.46%  "

# Row 6 - now Solana (was BNB)
$ws.Range("B6").Value = "Solana"
$ws.Range("C6").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D6").Value = "134.84"
$ws.Range("E6").Value = "  +2.42%  "

# Row 7 - LidoStakedEther
$ws.Range("D7").Value = "3.862.28"
$ws.Range("E7").Value = "  +8.37%  "

# Row 8 - XRP
$ws.Range("D8").Value = "0.621"
$ws.Range("E8").Value = "  -5.95%  "

# Row 9 - USDC
$ws.Range("E9").Value = "  +0.41%  "

# Row 10 - Cardano
$ws.Range("D10").Value = "0.731"
$ws.Range("E10").Value = "  -6.79%  "

# Row 11 - Dogecoin
$ws.Range("D11").Value = "0.164"
$ws.Range("E11").Value = "  -5.12%  "

# Row 12 - ShibaInu
$ws.Range("D12").Value = "0.0000314"
$ws.Range("E12").Value = "  +6.82%  "

# Row 13 - Avalanche
$ws.Range("D13").Value = "41.98"
$ws.Range("E13").Value = "  -3.47%  "

# Row 14 - now WrappedliquidstakedEther2.0 (was Polkadot)
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "4.361.00"
$ws.Range("E14").Value = "  +5.70%  "

# Row 15 - now Polkadot (was WrappedliquidstakedEther2.0)
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "9.91"
$ws.Range("E15").Value = "  -2.22%  "

# Row 16 - TRON
$ws.Range("D16").Value = "0.139"
$ws.Range("E16").Value = "  -1.10%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "3.772.88"
$ws.Range("E17").Value = "  +7.15%  "

# Row 18 - Chainlink
$ws.Range("D18").Value = "19.99"
$ws.Range("E18").Value = "  -2.40%  "

# Row 19 - Uniswap
$ws.Range("D19").Value = "12.95"
$ws.Range("E19").Value = "  +1.82%  "

# Row 20 - Polygon
$ws.Range("D20").Value = "1.08"
$ws.Range("E20").Value = "  -4.13%  "

# Row 21 - WrappedBTC
$ws.Range("D21").Value = "66.214.10"
$ws.Range("E21").Value = "  -0.45%  "

# Row 22 - BitcoinCash
$ws.Range("D22").Value = "418.08"
$ws.Range("E22").Value = "  -7.01%  "

# Row 23 - InternetComputer(DFINITY)
$ws.Range("D23").Value = "14.45"
$ws.Range("E23").Value = "  +9.94%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "86.18"
$ws.Range("E24").Value = "  -4.84%  "

# Row 25 - ImmutableX
$ws.Range("D25").Value = "2.97"
$ws.Range("E25").Value = "  -8.41%  "

# Row 26 - EthereumClassic
$ws.Range("D26").Value = "36.07"
$ws.Range("E26").Value = "  +4.98%  "

# Row 27 - PancakeSwap
$ws.Range("D27").Value = "3.22"
$ws.Range("E27").Value = "  -4.63%  "

# Row 28 - Filecoin
$ws.Range("D28").Value = "9.43"
$ws.Range("E28").Value = "  -6.31%  "

# Row 29 - LEO
$ws.Range("D29").Value = "5.15"
$ws.Range("E29").Value = "  +6.55%  "

# Row 30 - Cosmos
$ws.Range("D30").Value = "12.29"
$ws.Range("E30").Value = "  -1.72%  "

# Row 31 - Hedera
$ws.Range("D31").Value = "0.118"
$ws.Range("E31").Value = "  -1.12%  "

# Row 32 - Toncoin
$ws.Range("D32").Value = "2.69"
$ws.Range("E32").Value = "  -3.63%  "

# Row 33 - RenderToken
$ws.Range("D33").Value = "6.89"
$ws.Range("E33").Value = "  -6.34%  "

# Row 34 - now Kaspa (was InjectiveProtocol)
$ws.Range("B34").Value = "Kaspa"
$ws.Range("C34").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D34").Value = "0.158"
$ws.Range("E34").Value = "  -2.29%  "

# Row 35 - now InjectiveProtocol (was Kaspa)
$ws.Range("B35").Value = "InjectiveProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D35").Value = "40.36"
$ws.Range("E35").Value = "  +3.41%  "

# Row 36 - OKB
$ws.Range("D36").Value = "55.94"
$ws.Range("E36").Value = "  -2.72%  "

# Row 37 - Dai
$ws.Range("D37").Value = "0.997"
$ws.Range("E37").Value = "  -0.22%  "

# Row 38 - VeChain
$ws.Range("D38").Value = "0.0464"
$ws.Range("E38").Value = "  -7.88%  "

# Row 39 - ThetaToken
$ws.Range("D39").Value = "2.89"
$ws.Range("E39").Value = "  +23.53%  "

# Row 40 - now FirstDigitalUSD (was Stellar)
$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.31%  "

# Row 41 - now Stellar (was PEPE)
$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").Value = "0.139"
$ws.Range("E41").Value = "  -6.76%  "

# Row 42 - now EnergySwap (was FirstDigitalUSD)
$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").Value = "27.80"
$ws.Range("E42").Value = "  +27.07%  "

# Row 43 - now Monero (was LidoDAOToken)
$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D43").Value = "145.18"
$ws.Range("E43").Value = "  -2.28%  "

# Row 44 - now LidoDAOToken (was ARBITRUM)
$ws.Range("B44").Value = "LidoDAOToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D44").Value = "3.29"
$ws.Range("E44").Value = "  +0.64%  "

# Row 45 - ApeXProtocol
$ws.Range("D45").Value = "3.06"
$ws.Range("E45").Value = "  +19.47%  "

# Row 46 - now PEPE (was Monero)
$ws.Range("B46").Value = "PEPE"
$ws.Range("C46").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D46").Value = "0.0`u{2083}0618"
$ws.Range("E46").Value = "  -23.58%  "

# Row 47 - now ARBITRUM (was NEARProtocol)
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").Value = "2.03"
$ws.Range("E47").Value = "  +1.21%  "

# Row 48 - now NEARProtocol (was Stacks)
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "4.20"
$ws.Range("E48").Value = "  -5.71%  "

# Row 49 - now Stacks (was EnergySwap)
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").Value = "2.77"
$ws.Range("E49").Value = "  -8.95%  "

# Row 50 - WEMIXToken
$ws.Range("D50").Value = "2.54"
$ws.Range("E50").Value = "  -8.94%  "

# Row 51 - TheGraph
$ws.Range("D51").Value = "0.287"
$ws.Range("E51").Value = "  -7.50%  "
